$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) — update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1259
$ws1.Range("F4").Value = 2757
$ws1.Range("F5").Value = 246

# Sheet "全部类型" (all types) — same rows of data repeated, update accordingly
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1259
$ws4.Range("F6").Value = 2757
$ws4.Range("F8").Value = 246
